# Apply the benchmark-stats corrections to the single-column results table.
# Each change rewrites the text of one table cell (row, column 1) in place,
# which replaces whatever run(s)/text currently live in that cell's
# paragraph with a single run carrying the new value (inheriting the
# existing run formatting of the cell).

$d = $word.ActiveDocument
$t = $d.Tables(1)

$changes = @{
    1  = "0M"
    2  = "0M"
    3  = "0M"
    4  = "110"
    6  = "0.00010"
    7  = "0.00005"
    8  = "0.00002"
    9  = "0.00004"
    10 = "0.00004"
    11 = "0.00010"
    12 = "0.00477"
    44 = "100"
    45 = "0"
    46 = "170"
}

foreach ($row in $changes.Keys) {
    $cell = $t.Cell($row, 1)
    $cell.Range.Text = $changes[$row]
}
